# Refresh the cryptos list snapshot (prices / 1h volume %, plus the
# occasional rank swap between two coins with an identical price/volume)
# as scraped by the GitHub Actions job.
#
# Note: Range.Value= auto-coerces numeric-looking literals (e.g. "588.29")
# into real numbers, but column D in this sheet is always stored as text
# (even the original, unedited values are inlineStr). A leading apostrophe
# forces Excel to keep such literals as text (quotePrefix) without the
# apostrophe itself becoming part of the stored value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.009.43'
$ws.Range('E2').Value = '  -1.19%  '

$ws.Range('D3').Value = '2.592.35'
$ws.Range('E3').Value = '  -1.17%  '

$ws.Range('E4').Value = '  -0.14%  '

$ws.Range('D5').Value = '''588.29'
$ws.Range('E5').Value = '  -2.15%  '

$ws.Range('D6').Value = '''148.75'
$ws.Range('E6').Value = '  -3.99%  '

$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('D8').Value = '''0.545'
$ws.Range('E8').Value = '  -0.41%  '

$ws.Range('D9').Value = '2.590.51'
$ws.Range('E9').Value = '  -1.22%  '

$ws.Range('D10').Value = '''0.122'
$ws.Range('E10').Value = '  -3.89%  '

$ws.Range('E11').Value = '  -0.20%  '

$ws.Range('D12').Value = '''5.15'
$ws.Range('E12').Value = '  -2.10%  '

$ws.Range('D13').Value = '''0.342'
$ws.Range('E13').Value = '  -3.59%  '

$ws.Range('D14').Value = '''27.06'
$ws.Range('E14').Value = '  -2.34%  '

$ws.Range('D15').Value = '3.062.41'
$ws.Range('E15').Value = '  -1.32%  '

$ws.Range('E16').Value = '  -5.22%  '

$ws.Range('D17').Value = '66.981.33'
$ws.Range('E17').Value = '  -1.24%  '

$ws.Range('D18').Value = '2.590.25'
$ws.Range('E18').Value = '  -1.39%  '

$ws.Range('D19').Value = '''363.96'
$ws.Range('E19').Value = '  -0.71%  '

$ws.Range('D20').Value = '''10.91'
$ws.Range('E20').Value = '  -2.55%  '

$ws.Range('D21').Value = '''7.29'
$ws.Range('E21').Value = '  -4.65%  '

$ws.Range('E22').Value = '  -0.52%  '

$ws.Range('D23').Value = '''4.80'
$ws.Range('E23').Value = '  -2.68%  '

$ws.Range('D24').Value = '''2.02'
$ws.Range('E24').Value = '  -0.99%  '

$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '''72.05'
$ws.Range('E25').Value = '  +8.76%  '

$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  +0.04%  '

$ws.Range('D27').Value = '''9.88'
$ws.Range('E27').Value = '  +0.17%  '

$ws.Range('E29').Value = '  +0.03%  '

$ws.Range('D30').Value = '''576.67'
$ws.Range('E30').Value = '  +0.00%  '

$ws.Range('D31').Value = '0.0₃0970'
$ws.Range('E31').Value = '  -6.58%  '

$ws.Range('E32').Value = '  -4.82%  '

$ws.Range('D33').Value = '''7.57'
$ws.Range('E33').Value = '  -4.29%  '

$ws.Range('E34').Value = '  -3.83%  '

$ws.Range('E35').Value = '  -0.03%  '

$ws.Range('D36').Value = '''0.124'
$ws.Range('E36').Value = '  -6.47%  '

$ws.Range('D37').Value = '''1.48'
$ws.Range('E37').Value = '  -3.02%  '

$ws.Range('D38').Value = '''155.83'
$ws.Range('E38').Value = '  -1.34%  '

$ws.Range('D39').Value = '''18.87'
$ws.Range('E39').Value = '  -2.51%  '

$ws.Range('E40').Value = '  -1.60%  '

$ws.Range('D41').Value = '''1.83'
$ws.Range('E41').Value = '  -1.23%  '

$ws.Range('D42').Value = '''5.14'
$ws.Range('E42').Value = '  -3.70%  '

$ws.Range('D43').Value = '''16.78'
$ws.Range('E43').Value = '  +2.17%  '

$ws.Range('D44').Value = '''2.47'
$ws.Range('E44').Value = '  -3.94%  '

$ws.Range('E45').Value = '  -0.09%  '

$ws.Range('D46').Value = '''153.02'
$ws.Range('E46').Value = '  -2.56%  '

$ws.Range('D47').Value = '0.0₆0280'
$ws.Range('E47').Value = '  -1.92%  '

$ws.Range('E48').Value = '  -1.49%  '

$ws.Range('B49').Value = 'Optimism'
$ws.Range('C49').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D49').Value = '''1.67'
$ws.Range('E49').Value = '  -2.99%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '''0.0776'
$ws.Range('E50').Value = '  -1.74%  '

$ws.Range('D51').Value = '''21.18'
$ws.Range('E51').Value = '  +0.65%  '
